# mediclean_bd.xlsx bugfix: actualizar fechas de retiro del cliente,
# cerrar la ruta actual (moverla a rutas_bd como REALIZADO) y
# registrar el cierre en rutas_registros.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Hoja "clientes": adelantar la fecha de ultimo/proximo retiro
#    de Isaias Beroiza Mora (fila 3).
# ---------------------------------------------------------------
$clientes = $wb.Worksheets.Item("clientes")
$clientes.Range("I3").NumberFormat = "@"
$clientes.Range("I3").Value = "2024-08-16"
$clientes.Range("K3").NumberFormat = "@"
$clientes.Range("K3").Value = "2024-12-16"

# ---------------------------------------------------------------
# 2) Hoja "ruta_actual": se vacia (la ruta en curso termino), se
#    borran el encabezado (B1/C1) y la fila de datos (fila 3).
# ---------------------------------------------------------------
$rutaActual = $wb.Worksheets.Item("ruta_actual")
$rutaActual.Range("B1").ClearContents()
$rutaActual.Range("C1").ClearContents()
$rutaActual.Range("A3:L3").ClearContents()

# ---------------------------------------------------------------
# 3) Hoja "rutas_bd": se agrega como fila 8 el retiro que se acaba
#    de realizar (el mismo que tenia ruta_actual), marcado REALIZADO.
# ---------------------------------------------------------------
$rutasBd = $wb.Worksheets.Item("rutas_bd")
$rutasBd.Range("A8:L8").HorizontalAlignment = 1

$rutasBd.Range("A8").NumberFormat = "@"
$rutasBd.Range("A8").Value = "20240816"
$rutasBd.Range("B8").Value = 1
$rutasBd.Range("C8").Value = "16.742.249-7"
$rutasBd.Range("D8").Value = "Isaias Beroiza Mora"
$rutasBd.Range("E8").Value = "colaco sn km3 parcela 9"
$rutasBd.Range("F8").Value = "Calbuco"
$rutasBd.Range("G8").NumberFormat = "@"
$rutasBd.Range("G8").Value = "88809703"
$rutasBd.Range("H8").Value = "por buscar"
$rutasBd.Range("I8").Value = "ok"
$rutasBd.Range("J8").NumberFormat = "@"
$rutasBd.Range("J8").Value = "20240813"
$rutasBd.Range("K8").Value = "REALIZADO"
$rutasBd.Range("L8").Value = "ok"

# ---------------------------------------------------------------
# 4) Hoja "rutas_registros": se registra el cierre de "ruta TEST"
#    del 2024-08-16 (fila 5).
# ---------------------------------------------------------------
$rutasRegistros = $wb.Worksheets.Item("rutas_registros")
$rutasRegistros.Range("C5").Value = 1
$rutasRegistros.Range("E5").Value = "2024-08-16T15:37:51.807710 RUTA FINALIZADA"
